$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the extra rows (old rows 11-13) first so remaining rows shift up correctly
$ws.Rows("11:13").Delete()

# Row 3: Barcelona / Spain / Location
$ws.Range("A3").Value = "In wich country is Barcelona located?"
$ws.Range("B3").Value = "Spain"
$ws.Range("C3").Value = "Location"

# Row 4: Italian speakers / Italy / Location
$ws.Range("A4").Value = "Where do most people speak italian?"
$ws.Range("B4").Value = "Italy"
$ws.Range("C4").Value = "Location"

# Row 5: F1 World Champion 2022 / Max Verstappen / Person
$ws.Range("A5").Value = "Who was the F1 World Champion in 2022?"
$ws.Range("B5").Value = "Max Verstappen"
$ws.Range("C5").Value = "Person"

# Row 6: Mayor of Innsbruck / Georg Willi / Person
$ws.Range("A6").Value = "Who is the mayor of Innsbruck? "
$ws.Range("B6").Value = "Georg Willi"
$ws.Range("C6").Value = "Person"

# Row 7: Facebook founder / Marc Zuckerberg / Person
$ws.Range("A7").Value = "Who founded Facebook?"
$ws.Range("B7").Value = "Marc Zuckerberg"
$ws.Range("C7").Value = "Person"

# Row 8: Schumacher 1st title / 1994 / Year
$ws.Range("A8").Value = "When did Miachel Schumacher win his first F1 World Drivers Title?"
$ws.Range("B8").Value = 1994
$ws.Range("C8").Value = "Year"

# Row 9: Schumacher 3rd title / 2000 / Year
$ws.Range("A9").Value = "When did Miachel Schumacher win his 3rd F1 World Drivers Title?"
$ws.Range("B9").Value = 2000
$ws.Range("C9").Value = "Year"

# Row 10: Chelsea Champions League / 2021 / Year
$ws.Range("A10").Value = "When has Chelsea last won the Champions League?"
$ws.Range("B10").Value = 2021
$ws.Range("C10").Value = "Year"

# Update the selected cell shown in the saved view
$ws.Range("D20").Select()
